# Update EPEX Spot prices workbook:
#  - "Prix Spot" sheet: add a new day column CM ("12-sep") with its 24 hourly values.
#  - "Gaz" sheet: append a new row (88) for 2025-09-10.
#  - "CO2" sheet: append a new row (88) for 2025-09-10.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a literal text value into a cell without letting Excel's
# automatic type inference turn a date-looking string into a date serial
# number, and without leaving a stray NumberFormat style behind on the cell.
# ---------------------------------------------------------------------------
function Set-LiteralText {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# ---------------------------------------------------------------------------
# 1) "Prix Spot" sheet — add column CM ("12-sep")
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Header cell: copy the formatting of the previous header (CL1) so the new
# header matches the existing bold/bordered/centered style, then overwrite
# the value with the new date label.
$prevHeader = $wsPrix.Range("CL1")
$newHeader = $wsPrix.Range("CM1")
$prevHeader.Copy($newHeader)
$newHeader.Value = "12-sep"

# Hourly values for the new day (rows 2-25).
$prixValues = @{
    2  = 0
    3  = 0.01
    4  = 0.44
    5  = 0.01
    6  = 0.01
    7  = 3
    8  = 11.58
    9  = 21.04
    10 = 41.82
    11 = 26.68
    12 = 1.85
    13 = 0
    14 = 1.77
    15 = 0
    16 = 0
    17 = 0
    18 = 3.52
    19 = 19.68
    20 = 39.92
    21 = 61.68
    22 = 66.34
    23 = 46
    24 = 49.97
    25 = 41.2
}

foreach ($row in 2..25) {
    $wsPrix.Cells.Item($row, 91).Value = $prixValues[$row]
}

# ---------------------------------------------------------------------------
# 2) "Gaz" sheet — append row 88 for 2025-09-10
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
Set-LiteralText $wsGaz.Range("A88") "2025-09-10"
$wsGaz.Range("B88").Value = 32.6

# ---------------------------------------------------------------------------
# 3) "CO2" sheet — append row 88 for 2025-09-10
# ---------------------------------------------------------------------------
$wsCO2 = $wb.Worksheets.Item("CO2")
Set-LiteralText $wsCO2.Range("A88") "2025-09-10"
$wsCO2.Range("B88").Value = 76.2
